# İstisna Yönetimi ile yanlış giriş önleme
# Apply corrected/localized field values to the profile document using
# Find & Replace, guarded with try/catch so that a failed replacement
# (e.g. text not found) does not abort the whole script.

$d = $word.ActiveDocument

function Replace-Text($Old, $New) {
    try {
        $found = $d.Content.Find.Execute(
            $Old,   # FindText
            $true,  # MatchCase
            $true,  # MatchWholeWord
            $false, # MatchWildcards
            $false, # MatchSoundsLike
            $false, # MatchAllWordForms
            $true,  # Forward
            1,      # Wrap (wdFindContinue)
            $false, # Format
            $New,   # ReplaceWith
            2       # Replace (wdReplaceAll)
        )
        if (-not $found) {
            Write-Host "Warning: text not found -> '$Old'"
        }
    }
    catch {
        Write-Host "Error replacing '$Old': $_"
    }
}

Replace-Text "alex whooper" "Alex DeSouza"
Replace-Text "Gender: male" "Cinsiyet: Fenerbahçe"
Replace-Text "Date of birth: 11.11.2022" "Doğum tarihi: 11.11.1905"
Replace-Text "Marital status: married" "Medeni durumu: Bekar"
Replace-Text "Country: us" "Ülkesi: TR"
Replace-Text "Military status: done" "Askerlik durumu: Yapıldı"
Replace-Text "License type: b" "Ehliyet türü: B"
